$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9535.333000000001
$ws.Range("I32").Value = 10964.2
$ws.Range("J32").Value = 7749.25
$ws.Range("K32").Value = 10964.2
$ws.Range("L32").Value = 7749.25
$ws.Range("M32").Value = -10638.2
$ws.Range("N32").Value = -8401.25
$ws.Range("H33").Value = 169.52
$ws.Range("I33").Value = 178.60869
$ws.Range("K33").Value = 178.60869
$ws.Range("M33").Value = 50.39131
$ws.Range("H43").Value = 4013.0356
$ws.Range("I43").Value = 2123
$ws.Range("K43").Value = 2123
$ws.Range("M43").Value = -2054
$ws.Range("H61").Value = 1774.2222
$ws.Range("I61").Value = 1774.2222
$ws.Range("K61").Value = 5322.6666
$ws.Range("M61").Value = -5150.6666
$ws.Range("H74").Value = 8018.9443
$ws.Range("I74").Value = 6048.909
$ws.Range("K74").Value = 6048.909
$ws.Range("M74").Value = -5112.909
$ws.Range("H77").Value = 8018.9443
$ws.Range("I77").Value = 6048.909
$ws.Range("K77").Value = 30244.545
$ws.Range("M77").Value = -25564.545
$ws.Range("H132").Value = 1518.3684
$ws.Range("I132").Value = 1456.6
$ws.Range("K132").Value = 4369.799999999999
$ws.Range("M132").Value = -1839.799999999999
$ws.Range("H138").Value = 3883.35
$ws.Range("I138").Value = 5743.25
$ws.Range("J138").Value = 3418.375
$ws.Range("K138").Value = 17229.75
$ws.Range("L138").Value = 10255.125
$ws.Range("M138").Value = -12089.75
$ws.Range("N138").Value = -20535.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1905.2034
$ws.Range("I32").Value = 1504.4182
$ws.Range("K32").Value = 1504.4182
$ws.Range("M32").Value = -1217.4182
$ws.Range("H45").Value = 3814.2942
$ws.Range("I45").Value = 1994.9286
$ws.Range("K45").Value = 1994.9286
$ws.Range("M45").Value = -1617.9286
$ws.Range("H122").Value = 111115176
$ws.Range("I122").Value = 1799.5
$ws.Range("K122").Value = 5398.5
$ws.Range("M122").Value = -2948.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3588
$ws.Range("J16").Value = 4563.875
$ws.Range("L16").Value = 4563.875
$ws.Range("N16").Value = -5137.875
$ws.Range("H22").Value = 1899.9166
$ws.Range("I22").Value = 416.66666
$ws.Range("J22").Value = 3383.1667
$ws.Range("K22").Value = 416.66666
$ws.Range("L22").Value = 3383.1667
$ws.Range("M22").Value = -66.66665999999998
$ws.Range("N22").Value = -4083.1667
$ws.Range("H113").Value = 3588
$ws.Range("J113").Value = 4563.875
$ws.Range("L113").Value = 4563.875
$ws.Range("N113").Value = -8903.875
$ws.Range("H132").Value = 2201.3
$ws.Range("I132").Value = 1118.9412
$ws.Range("J132").Value = 8334.666999999999
$ws.Range("K132").Value = 3356.8236
$ws.Range("L132").Value = 25004.001
$ws.Range("M132").Value = -826.8235999999997
$ws.Range("N132").Value = -30064.001
$ws.Range("H134").Value = 3871.5
$ws.Range("I134").Value = 2178.7
$ws.Range("K134").Value = 6536.099999999999
$ws.Range("M134").Value = -4001.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6513.9
$ws.Range("I5").Value = 2716.6667
$ws.Range("K5").Value = 8150.000100000001
$ws.Range("M5").Value = -8038.000100000001
$ws.Range("H12").Value = 22.6875
$ws.Range("I12").Value = 67.666664
$ws.Range("J12").Value = 12.307693
$ws.Range("K12").Value = 202.999992
$ws.Range("L12").Value = 36.923079
$ws.Range("M12").Value = -29.99999199999999
$ws.Range("N12").Value = -382.923079
$ws.Range("H128").Value = 201317.33
$ws.Range("I128").Value = 201317.33
$ws.Range("K128").Value = 603951.99
$ws.Range("M128").Value = -598971.99
$ws.Range("H135").Value = 6513.9
$ws.Range("I135").Value = 2716.6667
$ws.Range("K135").Value = 24450.0003
$ws.Range("M135").Value = -21915.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""
$ws.Range("H80").Value = 10675.5
$ws.Range("I80").Value = 8000
$ws.Range("K80").Value = 8000
$ws.Range("M80").Value = -7002
$ws.Range("H83").Value = 10675.5
$ws.Range("I83").Value = 8000
$ws.Range("K83").Value = 40000
$ws.Range("M83").Value = -35008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2435.8333
$ws.Range("I16").Value = 922.7
$ws.Range("J16").Value = 10001.5
$ws.Range("K16").Value = 922.7
$ws.Range("L16").Value = 10001.5
$ws.Range("M16").Value = -752.7
$ws.Range("N16").Value = -10341.5
$ws.Range("H46").Value = 2333
$ws.Range("J46").Value = 2999.5
$ws.Range("L46").Value = 2999.5
$ws.Range("N46").Value = -3375.5
$ws.Range("H97").Value = 8454.546
$ws.Range("J97").Value = 8523.809999999999
$ws.Range("L97").Value = 8523.809999999999
$ws.Range("N97").Value = -10505.81
$ws.Range("H100").Value = 11753
$ws.Range("I100").Value = 8529.615
$ws.Range("J100").Value = 15245
$ws.Range("K100").Value = 8529.615
$ws.Range("L100").Value = 15245
$ws.Range("M100").Value = -7988.615
$ws.Range("N100").Value = -16327
$ws.Range("H101").Value = 16574.75
$ws.Range("J101").Value = 16574.75
$ws.Range("L101").Value = 16574.75
$ws.Range("N101").Value = -23064.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 701
$ws.Range("I113").Value = 701.2857
$ws.Range("K113").Value = 2103.8571
$ws.Range("M113").Value = 66.14289999999983
$ws.Range("H122").Value = 5566
$ws.Range("I122").Value = 1575.2
$ws.Range("J122").Value = 15543
$ws.Range("K122").Value = 4725.6
$ws.Range("L122").Value = 46629
$ws.Range("M122").Value = -2275.6
$ws.Range("N122").Value = -51529
$ws.Range("H126").Value = 3928.4285
$ws.Range("I126").Value = 3849.5
$ws.Range("K126").Value = 11548.5
$ws.Range("M126").Value = -9078.5
$ws.Range("H132").Value = 6342.4116
$ws.Range("J132").Value = 14015.1
$ws.Range("L132").Value = 42045.3
$ws.Range("N132").Value = -47105.3
$ws.Range("H136").Value = 2429.2
$ws.Range("I136").Value = 1937.8948
$ws.Range("J136").Value = 3985
$ws.Range("K136").Value = 5813.6844
$ws.Range("L136").Value = 11955
$ws.Range("M136").Value = -3263.6844
$ws.Range("N136").Value = -17055
